$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "-"
$ws.Range("E2").Value = "-"

# Row 8
$ws.Range("D8").Value = "-"
$ws.Range("E8").Value = "-"

# Rows 11,12,14,15
$ws.Range("E11").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("E15").Value = "-"

# Row 18
$ws.Range("B18").Value = "-"
$ws.Range("F18").Value = "['ELM-2NA-CAM', -, -, -]"

# Row 19
$ws.Range("B19").Value = "-"
$ws.Range("D19").Value = "[-, -, -, 'ELM-2NA-CAM']"
$ws.Range("F19").Value = "[-, 'ELM-2NA-CAM', -, -]"

# Row 20
$ws.Range("B20").Value = "-"
$ws.Range("F20").Value = "[-, 'ELM-2NA-CAM', -, -]"

# Row 21
$ws.Range("B21").Value = "-"
$ws.Range("F21").Value = "-"
